# Add a "Save" column (column H) to the s_vals sheet, mirroring the header
# formatting already used by the other header cells (e.g. G1), and fill in
# the data values for the two existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, centered/top alignment, thin
# border) from the last existing header cell (G1) onto the new header cell.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and the data values for rows 2-3.
$ws.Range("H1").Value2 = "Save"
$ws.Range("H2").Value2 = 1
$ws.Range("H3").Value2 = 0
